# Apply the updated cryptocurrency price/volume figures to the worksheet.
# Numeric-looking "Price" strings (e.g. "64.70") are written with a leading
# apostrophe so Excel stores them as text (matching the source data, which
# keeps trailing zeros / non-standard grouping like "58.330.84") instead of
# silently re-parsing them into a Double and losing formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.330.84'
$ws.Range('E2').Value = '  -4.38%  '
$ws.Range('D3').Value = '2.646.60'
$ws.Range('E3').Value = '  -1.84%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '''521.83'
$ws.Range('E5').Value = '  -1.07%  '
$ws.Range('D6').Value = '''144.38'
$ws.Range('E6').Value = '  -0.34%  '
$ws.Range('D7').Value = '''0.999'
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('E8').Value = '  -2.02%  '
$ws.Range('D9').Value = '''6.68'
$ws.Range('E9').Value = '  +0.27%  '
$ws.Range('D10').Value = '''0.103'
$ws.Range('E10').Value = '  -3.13%  '
$ws.Range('E11').Value = '  -0.68%  '
$ws.Range('D12').Value = '''0.132'
$ws.Range('E12').Value = '  +1.59%  '
$ws.Range('D13').Value = '3.109.55'
$ws.Range('E13').Value = '  -2.07%  '
$ws.Range('D14').Value = '58.321.62'
$ws.Range('E14').Value = '  -4.28%  '
$ws.Range('D15').Value = '''20.89'
$ws.Range('E15').Value = '  -2.05%  '
$ws.Range('E16').Value = '  -1.52%  '
$ws.Range('D17').Value = '2.648.70'
$ws.Range('E17').Value = '  -2.28%  '
$ws.Range('D18').Value = '''338.61'
$ws.Range('E18').Value = '  -3.00%  '
$ws.Range('D19').Value = '''4.41'
$ws.Range('E19').Value = '  -2.56%  '
$ws.Range('D20').Value = '''10.47'
$ws.Range('E20').Value = '  -1.03%  '
$ws.Range('E21').Value = '  -1.00%  '
$ws.Range('E22').Value = '  +0.19%  '
$ws.Range('D23').Value = '''64.70'
$ws.Range('E23').Value = '  +1.55%  '
$ws.Range('D24').Value = '''0.425'
$ws.Range('E24').Value = '  +0.91%  '
$ws.Range('E25').Value = '  -1.31%  '
$ws.Range('E26').Value = '  +0.39%  '
$ws.Range('D27').Value = '0.0₃0797'
$ws.Range('E27').Value = '  -2.52%  '
$ws.Range('D28').Value = '''7.11'
$ws.Range('E28').Value = '  -2.83%  '
$ws.Range('E29').Value = '  -2.10%  '
$ws.Range('E30').Value = '  +0.00%  '
$ws.Range('E31').Value = '  -1.04%  '
$ws.Range('D32').Value = '''152.63'
$ws.Range('E32').Value = '  +1.62%  '
$ws.Range('D33').Value = '''18.83'
$ws.Range('E33').Value = '  -1.62%  '
$ws.Range('D34').Value = '''4.14'
$ws.Range('E34').Value = '  -2.95%  '
$ws.Range('D35').Value = '''0.915'
$ws.Range('E35').Value = '  -3.82%  '
$ws.Range('D37').Value = '''0.859'
$ws.Range('E37').Value = '  -2.51%  '
$ws.Range('D38').Value = '''36.79'
$ws.Range('E39').Value = '  -5.79%  '
$ws.Range('D40').Value = '''3.65'
$ws.Range('E40').Value = '  -0.71%  '
$ws.Range('E41').Value = '  +0.22%  '
$ws.Range('D42').Value = '''0.606'
$ws.Range('E42').Value = '  -0.90%  '
$ws.Range('D43').Value = '''0.0968'
$ws.Range('E43').Value = '  -2.66%  '
$ws.Range('D44').Value = '''270.08'
$ws.Range('E44').Value = '  -5.54%  '
$ws.Range('D45').Value = '''19.40'
$ws.Range('E45').Value = '  -3.11%  '
$ws.Range('D46').Value = '''0.0537'
$ws.Range('E46').Value = '  -0.70%  '
$ws.Range('E47').Value = '  +1.42%  '
$ws.Range('D48').Value = '2.040.87'
$ws.Range('E48').Value = '  -5.29%  '
$ws.Range('E49').Value = '  -3.27%  '
$ws.Range('D50').Value = '''4.66'
$ws.Range('E50').Value = '  -3.15%  '
$ws.Range('D51').Value = '''18.34'
$ws.Range('E51').Value = '  -4.14%  '
